# Add files via upload
# - Populate Minute/Second/Rep (columns O,P,Q) for the 6th workout on the
#   "ScoreM" sheet, rows 2-11.
# - Switch the active sheet/tab back to "ScoreM" (was "ScoreF"), updating the
#   selection on each sheet accordingly.

$wb = $excel.ActiveWorkbook

$wsScoreM = $wb.Worksheets.Item("ScoreM")

# row, Minute(O), Second(P), Rep(Q)
$values = @(
    , @(2,  8, 0, 180)
    , @(3,  8, 0, 220)
    , @(4,  8, 0, 144)
    , @(5,  8, 0, 174)
    , @(6,  8, 0, 238)
    , @(7,  8, 0, 258)
    , @(8,  8, 0, 259)
    , @(9,  8, 0, 240)
    , @(10, 8, 0, 238)
    , @(11, 8, 0, 253)
)

foreach ($entry in $values) {
    $row = $entry[0]
    $wsScoreM.Range("O$row").Value = $entry[1]
    $wsScoreM.Range("P$row").Value = $entry[2]
    $wsScoreM.Range("Q$row").Value = $entry[3]
}

# Make "ScoreM" the active sheet/tab again and restore its selection.
$wsScoreM.Activate()
$wsScoreM.Range("A5").Select()
